$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add a new "Test case 7" block (rows 61-69) for the calculateVolume
# function, mirroring the structure of the preceding "Test case 6"
# block (rows 51-59).
# ------------------------------------------------------------------

# 1) Create the merged regions first (while the target cells are still
#    blank) so the merge operation does not re-split any border styles
#    that get pasted onto them afterwards.
$ws.Range("A61:F61").Merge()
$ws.Range("C63:C64").Merge()
$ws.Range("D63:D64").Merge()
$ws.Range("A68:A69").Merge()
$ws.Range("B68:B69").Merge()
$ws.Range("D68:D69").Merge()
$ws.Range("E68:E69").Merge()
$ws.Range("F68:F69").Merge()

# 2) Copy the cell formatting from the analogous cells of the
#    "Test case 6" block (rows 51-59) onto the new block (rows 61-69).
$ws.Range("A51").Copy()
$ws.Range("A61").PasteSpecial(-4122)
$ws.Range("B51:F51").Copy()
$ws.Range("B61:F61").PasteSpecial(-4122)

$ws.Range("A53").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("B53").Copy()
$ws.Range("B63").PasteSpecial(-4122)
$ws.Range("C53").Copy()
$ws.Range("C63").PasteSpecial(-4122)
$ws.Range("D53").Copy()
$ws.Range("D63").PasteSpecial(-4122)

$ws.Range("A54").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("B54").Copy()
$ws.Range("B64").PasteSpecial(-4122)
$ws.Range("C54").Copy()
$ws.Range("C64").PasteSpecial(-4122)
$ws.Range("D54").Copy()
$ws.Range("D64").PasteSpecial(-4122)

$ws.Range("A55").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("B55").Copy()
$ws.Range("B65").PasteSpecial(-4122)
$ws.Range("C55").Copy()
$ws.Range("C65").PasteSpecial(-4122)
$ws.Range("D55").Copy()
$ws.Range("D65").PasteSpecial(-4122)

$ws.Range("A57:F57").Copy()
$ws.Range("A67:F67").PasteSpecial(-4122)

$ws.Range("A58").Copy()
$ws.Range("A68").PasteSpecial(-4122)
$ws.Range("B58").Copy()
$ws.Range("B68").PasteSpecial(-4122)
$ws.Range("C58").Copy()
$ws.Range("C68").PasteSpecial(-4122)
$ws.Range("D58").Copy()
$ws.Range("D68").PasteSpecial(-4122)
$ws.Range("E58").Copy()
$ws.Range("E68").PasteSpecial(-4122)
$ws.Range("F58").Copy()
$ws.Range("F68").PasteSpecial(-4122)

$ws.Range("A59").Copy()
$ws.Range("A69").PasteSpecial(-4122)
$ws.Range("B59").Copy()
$ws.Range("B69").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("C69").PasteSpecial(-4122)
$ws.Range("D59").Copy()
$ws.Range("D69").PasteSpecial(-4122)
$ws.Range("E59").Copy()
$ws.Range("E69").PasteSpecial(-4122)
$ws.Range("F59").Copy()
$ws.Range("F69").PasteSpecial(-4122)

# Match the auto-fit row height used by the other section-title rows.
$ws.Rows(61).RowHeight = $ws.Rows(51).RowHeight

# 3) Fill in the text/values for the new block.
$ws.Range("A61").Value = "Test case 7"

$ws.Range("A63").Value = "ID"
$ws.Range("B63").Value = "UT_002"
$ws.Range("C63").Value = "Created by"
$ws.Range("D63").Value = "Trani Tranev"

$ws.Range("A64").Value = "Name"
$ws.Range("B64").Value = "UnitTesting"

$ws.Range("A65").Value = "Description"
$ws.Range("B65").Value = "Test the function: calculateVolume"
$ws.Range("C65").Value = "Created on"
$ws.Range("D65").Value = 44534

$ws.Range("A67").Value = "#"
$ws.Range("B67").Value = "Description"
$ws.Range("C67").Value = "Test Data"
$ws.Range("D67").Value = "Expectations"
$ws.Range("E67").Value = "Actual Result"
$ws.Range("F67").Value = "Status"

$ws.Range("A68").Value = 1
$ws.Range("B68").Value = "Testing calculateVolume function"
$ws.Range("C68").Value = 20
$ws.Range("D68").Value = 3768
$ws.Range("E68").Value = 3768
$ws.Range("F68").Value = "PASS"

$ws.Range("C69").Value = 3
